# Add notes on linear algebra: give the "notes" cell for week 1 / Warm-up I
# its own distinct text ("warm-up-1") instead of reusing the "w1p1" slides
# link, and clear out the now-unused slide/assignment columns for the
# remaining rows (they were placeholder duplicates of the "Final" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E2 ("notes" column, week 1 / Warm-up I row) gets its own unique text.
$ws.Range("E2").Value = "warm-up-1"

# Row 3 no longer needs slide links (D3:E3), and rows 3-6 no longer need
# the assigned/due link+title placeholder columns (G:J).
$ws.Range("D3:E3").Clear()
$ws.Range("G3:J6").Clear()

# Update the view state left behind by the edit: zoomed out a bit and the
# selection moved off the old D3:E3 range onto I16.
$excel.ActiveWindow.Zoom = 70
$ws.Range("I16").Select()
